# Auto-generated COM-interop script.
# Re-shuffles the existing 46 data rows (2-47) of "Avverkningsanmälningar"
# into their new positions (rows keep their own formulas / number formats,
# only their row number and the "Forandrad" (C) date change), then appends
# three brand-new rows (45, 48, 50) with literal data. Net effect: the used
# range grows from A1:Z47 to A1:Z50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ColCount = 26
$NewDate = 46079

# ---------------------------------------------------------------------
# 1) Snapshot every cell (Formula / NumberFormat / WrapText) of the
#    currently-populated rows 2..47 BEFORE any writes happen, since the
#    target layout re-uses/relocates this same data (e.g. rows 5 and 6
#    swap places), so source and destination ranges overlap.
# ---------------------------------------------------------------------
$snapshot = @{}
for ($r = 2; $r -le 47; $r++) {
    $rowData = @()
    for ($c = 1; $c -le $ColCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $rowData += [PSCustomObject]@{
            Formula      = $cell.Formula
            NumberFormat = $cell.NumberFormat
            WrapText     = $cell.WrapText
        }
    }
    $snapshot[$r] = $rowData
}

# ---------------------------------------------------------------------
# 2) new-row -> old-row mapping: the whole row (every column, every
#    formula / style) moves from the old row number to the new one.
# ---------------------------------------------------------------------
$rowMap = @{2=2; 3=3; 4=4; 5=6; 6=5; 7=7; 8=8; 9=9; 10=10; 11=11; 12=12; 13=13; 14=14; 15=15; 16=16; 17=33; 18=43; 19=17; 20=36; 21=39; 22=30; 23=35; 24=29; 25=18; 26=24; 27=23; 28=19; 29=21; 30=20; 31=45; 32=22; 33=47; 34=25; 35=31; 36=46; 37=27; 38=26; 39=28; 40=32; 41=44; 42=34; 43=37; 44=38; 46=42; 47=41; 49=40}

function Write-MappedRow {
    param($destRow, $srcRow)
    $src = $snapshot[$srcRow]
    for ($c = 1; $c -le $ColCount; $c++) {
        $srcCell = $src[$c - 1]
        $dst = $ws.Cells.Item($destRow, $c)
        $dst.Formula = $srcCell.Formula
        if ($c -eq 2 -or $c -eq 3) {
            $dst.NumberFormat = $srcCell.NumberFormat
        }
        if ($c -eq 18) {
            $dst.WrapText = $srcCell.WrapText
        }
    }
    # Column C ("Forandrad") is bumped to the new date for every row.
    $cCell = $ws.Cells.Item($destRow, 3)
    $cCell.Formula = $NewDate
    $cCell.NumberFormat = "YYYY-MM-DD"
}

foreach ($destRow in ($rowMap.Keys | Sort-Object { [int]$_ })) {
    $srcRow = $rowMap[$destRow]
    Write-MappedRow $destRow $srcRow
}

# ---------------------------------------------------------------------
# 3) Brand-new rows with no predecessor (45, 48, 50): literal data.
# ---------------------------------------------------------------------
function Write-NewRow {
    param($Row, $Data, $SetHt)

    foreach ($colLetter in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")) {
        if ($Data.ContainsKey($colLetter)) {
            $colIndex = [int][char]$colLetter - [int][char]'A' + 1
            $cell = $ws.Cells.Item($Row, $colIndex)
            $cell.Formula = $Data[$colLetter]
            if ($colLetter -eq "B" -or $colLetter -eq "C") {
                $cell.NumberFormat = "YYYY-MM-DD"
            }
        }
    }
    # "Artnamn" (R, column 18) is always present (even if empty) with the
    # wrap-text style used throughout the table.
    $rCell = $ws.Cells.Item($Row, 18)
    $rCell.WrapText = $true

    if ($SetHt) {
        $ws.Rows.Item($Row).RowHeight = 15
    }
}

$newRow = @{}
$newRow['A'] = 'A 10654-2026'
$newRow['B'] = 46078.35840277778
$newRow['C'] = 46079
$newRow['D'] = 'VÄSTRA GÖTALANDS LÄN'
$newRow['E'] = 'MÖLNDAL'
$newRow['F'] = 'Kommuner'
$newRow['G'] = 2.4
$newRow['H'] = 0
$newRow['I'] = 0
$newRow['J'] = 0
$newRow['K'] = 0
$newRow['L'] = 0
$newRow['M'] = 0
$newRow['N'] = 0
$newRow['O'] = 0
$newRow['P'] = 0
$newRow['Q'] = 0
Write-NewRow 45 $newRow $true

$newRow = @{}
$newRow['A'] = 'A 10672-2026'
$newRow['B'] = 46078.40140046296
$newRow['C'] = 46079
$newRow['D'] = 'VÄSTRA GÖTALANDS LÄN'
$newRow['E'] = 'MÖLNDAL'
$newRow['F'] = 'Kommuner'
$newRow['G'] = 1.4
$newRow['H'] = 0
$newRow['I'] = 0
$newRow['J'] = 0
$newRow['K'] = 0
$newRow['L'] = 0
$newRow['M'] = 0
$newRow['N'] = 0
$newRow['O'] = 0
$newRow['P'] = 0
$newRow['Q'] = 0
Write-NewRow 48 $newRow $true

$newRow = @{}
$newRow['A'] = 'A 10661-2026'
$newRow['B'] = 46078.37611111111
$newRow['C'] = 46079
$newRow['D'] = 'VÄSTRA GÖTALANDS LÄN'
$newRow['E'] = 'MÖLNDAL'
$newRow['F'] = 'Kommuner'
$newRow['G'] = 1.9
$newRow['H'] = 0
$newRow['I'] = 0
$newRow['J'] = 0
$newRow['K'] = 0
$newRow['L'] = 0
$newRow['M'] = 0
$newRow['N'] = 0
$newRow['O'] = 0
$newRow['P'] = 0
$newRow['Q'] = 0
Write-NewRow 50 $newRow $false


# ---------------------------------------------------------------------
# 4) Make sure the sheet's used range is recognised as A1:Z50.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 26).Formula = $ws.Cells.Item(1, 26).Formula
$ws.Cells.Item(50, 1).Formula = $ws.Cells.Item(50, 1).Formula
